$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (the "R40" rule row) gets its "Rule" label (column B) renamed to "1".
# Format the cell as Text first so Excel stores the digit-only value as a
# literal string (shared string "1") instead of re-interpreting it as the
# number 1 -- matching how the workbook now has an extra shared string "1"
# feeding cell B11.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
